$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.909.92"
$ws.Range("E2").Value = "  -5.49%  "
$ws.Range("D3").Value = "3.738.71"
$ws.Range("E3").Value = "  -6.25%  "
$ws.Range("E4").Value = "  +0.25%  "
$ws.Range("D5").Value = "'573.80"
$ws.Range("E5").Value = "  -3.25%  "
$ws.Range("D6").Value = "'160.53"
$ws.Range("E6").Value = "  -0.49%  "
$ws.Range("D7").Value = "'0.645"
$ws.Range("E7").Value = "  -5.70%  "
$ws.Range("E8").Value = "  +0.40%  "
$ws.Range("D9").Value = "'0.719"
$ws.Range("E9").Value = "  -3.92%  "
$ws.Range("E10").Value = "  -0.09%  "
$ws.Range("D11").Value = "'50.93"
$ws.Range("E11").Value = "  -5.62%  "
$ws.Range("D12").Value = "'0.0000307"
$ws.Range("E12").Value = "  -3.56%  "
$ws.Range("D13").Value = "'10.78"
$ws.Range("E13").Value = "  -1.63%  "
$ws.Range("D14").Value = "4.337.54"
$ws.Range("E14").Value = "  -6.12%  "
$ws.Range("D15").Value = "3.758.83"
$ws.Range("E15").Value = "  -5.90%  "
$ws.Range("D16").Value = "'20.06"
$ws.Range("E16").Value = "  -1.49%  "
$ws.Range("E17").Value = "  -8.24%  "
$ws.Range("D18").Value = "'13.29"
$ws.Range("E18").Value = "  -6.12%  "
$ws.Range("E19").Value = "  -2.82%  "
$ws.Range("D20").Value = "68.691.90"
$ws.Range("E20").Value = "  -5.33%  "
$ws.Range("D21").Value = "'420.61"
$ws.Range("E21").Value = "  -2.91%  "
$ws.Range("D22").Value = "'4.56"
$ws.Range("E22").Value = "  -4.79%  "
$ws.Range("D23").Value = "'90.58"
$ws.Range("E23").Value = "  -5.84%  "
$ws.Range("D24").Value = "'3.16"
$ws.Range("E24").Value = "  -7.81%  "
$ws.Range("D25").Value = "'13.39"
$ws.Range("E25").Value = "  -5.60%  "
$ws.Range("D26").Value = "'10.80"
$ws.Range("E26").Value = "  -4.53%  "
$ws.Range("D27").Value = "'3.81"
$ws.Range("E27").Value = "  -13.27%  "
$ws.Range("D28").Value = "'5.88"
$ws.Range("E28").Value = "  -1.26%  "
$ws.Range("D29").Value = "'9.99"
$ws.Range("E29").Value = "  -4.83%  "
$ws.Range("D30").Value = "'33.75"
$ws.Range("E30").Value = "  -7.19%  "
$ws.Range("D31").Value = "'7.79"
$ws.Range("E31").Value = "  -1.42%  "
$ws.Range("D32").Value = "'13.01"
$ws.Range("E32").Value = "  -5.10%  "
$ws.Range("D33").Value = "'46.15"
$ws.Range("E33").Value = "  -5.51%  "
$ws.Range("E34").Value = "  -7.73%  "
$ws.Range("D35").Value = "'67.36"
$ws.Range("E35").Value = "  -4.31%  "
$ws.Range("D36").Value = "0.0₃0933"
$ws.Range("E36").Value = "  +5.83%  "
$ws.Range("D37").Value = "'611.98"
$ws.Range("E37").Value = "  -8.85%  "
$ws.Range("D38").Value = "'0.410"
$ws.Range("E38").Value = "  -6.08%  "
$ws.Range("D39").Value = "'0.999"
$ws.Range("E39").Value = "  +0.20%  "
$ws.Range("E40").Value = "  +0.23%  "
$ws.Range("E41").Value = "  -4.20%  "
$ws.Range("E42").Value = "  -6.43%  "
$ws.Range("D43").Value = "'3.06"
$ws.Range("E43").Value = "  +16.14%  "
$ws.Range("D44").Value = "'0.0452"
$ws.Range("E44").Value = "  -7.76%  "
$ws.Range("D45").Value = "'2.73"
$ws.Range("E45").Value = "  +4.91%  "
$ws.Range("D46").Value = "'9.55"
$ws.Range("E46").Value = "  -10.65%  "
$ws.Range("D47").Value = "'0.139"
$ws.Range("E47").Value = "  -7.12%  "
$ws.Range("E48").Value = "  -17.17%  "
$ws.Range("D49").Value = "2.778.19"
$ws.Range("E49").Value = "  -2.24%  "
$ws.Range("D50").Value = "'3.18"
$ws.Range("E50").Value = "  -8.04%  "
$ws.Range("D51").Value = "'0.000262"
$ws.Range("E51").Value = "  -2.07%  "
